$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "listOfCities"
$ws.Name = "listOfCities"

# Reset the active selection back to A1 (the sheet previously had A5 selected)
$ws.Range("A1").Select()
